$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (shifts the existing
# 2025-03-11 data row down to row 3).
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new data point
# (date 2025-04-28 => serial 45775, value 2.3%).
$ws.Range("A2").Value = 45775
$ws.Range("B2").Value = 0.023

# Copy number formatting (date / percent) from the row below, which
# still carries the original cell styles, onto the new row without
# disturbing its values.
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
